# Update the crypto price/volume snapshot (Price column D, Volume(1h) column E).
# All D/E cells are stored as text in the workbook; values in column D that look
# like plain numbers are prefixed with a leading apostrophe so Excel keeps
# storing them as text instead of auto-converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.549.31"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "1.826.12"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D5").Value = "'315.62"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.5118"
$ws.Range("E7").Value = "  -5.56%  "
$ws.Range("D8").Value = "'0.3957"
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("D9").Value = "'0.08196"
$ws.Range("E9").Value = "  +6.71%  "
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").Value = "'41.74"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "'21.20"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "'6.346"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").Value = "'7.546"
$ws.Range("D16").Value = "1.826.75"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "'0.00001128"
$ws.Range("E17").Value = "  +3.24%  "
$ws.Range("D18").Value = "'92.89"
$ws.Range("E18").Value = "  +3.10%  "
$ws.Range("D19").Value = "'0.06665"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").Value = "'17.84"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'6.100"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").Value = "28.584.90"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'11.44"
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("D25").Value = "'2.262"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").Value = "'21.46"
$ws.Range("E26").Value = "  +3.26%  "
$ws.Range("D27").Value = "'156.83"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").Value = "2.035.46"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "'2.412"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("D30").Value = "'127.49"
$ws.Range("E30").Value = "  +2.58%  "
$ws.Range("D31").Value = "'1.113"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").Value = "'5.774"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("E35").Value = "  -6.58%  "
$ws.Range("D36").Value = "'0.2233"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").Value = "'5.289"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").Value = "'0.02356"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").Value = "'8.823"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").Value = "'0.6332"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").Value = "'11.31"
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("D42").Value = "'1.181"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").Value = "'1.399"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").Value = "'13.58"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("D45").Value = "'0.5948"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("D46").Value = "'3.735"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").Value = "'125.47"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "'1.998"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").Value = "'1.196"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").Value = "'0.06932"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").Value = "'1.087"
$ws.Range("E51").Value = "  +4.65%  "
